# "semana 32 de 2025" - add the week-32 column (AI) to the weekly IRA
# extract: a new header label in AI1 and one case count per UPGD row
# that reported data for that week (rows with no week-32 record are
# left untouched, matching the sparse source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell AI1 holds the week number "32" as text (like the existing
# week-number headers D1:AH1, which are all text, not numbers). Prefixing
# with a single quote forces Excel to store the digits as text instead of
# coercing them to a number, while keeping the bold/centered header look.
$ws.Range("AI1").Value = "'32"

# Week-32 case counts, one per UPGD (row).
$weekCounts = [ordered]@{
    "AI2"  = 49
    "AI3"  = 73
    "AI5"  = 4
    "AI6"  = 61
    "AI7"  = 31
    "AI8"  = 30
    "AI9"  = 1
    "AI10" = 3
    "AI13" = 1
    "AI15" = 2
    "AI16" = 1
    "AI17" = 2
    "AI23" = 2
    "AI25" = 53
    "AI27" = 1
    "AI28" = 0
    "AI29" = 0
    "AI30" = 16
    "AI31" = 4
    "AI32" = 1
    "AI34" = 0
    "AI35" = 43
    "AI36" = 5
    "AI37" = 10
    "AI38" = 68
    "AI40" = 2
    "AI41" = 12
    "AI42" = 8
    "AI43" = 29
    "AI45" = 64
    "AI46" = 131
    "AI47" = 4
    "AI48" = 106
    "AI49" = 5
    "AI50" = 0
    "AI51" = 11
    "AI53" = 32
    "AI54" = 1
    "AI55" = 0
    "AI56" = 4
    "AI57" = 24
    "AI58" = 21
}

foreach ($cellRef in $weekCounts.Keys) {
    $ws.Range($cellRef).Value = $weekCounts[$cellRef]
}
